$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the N2 header cell from "Speed" to "SPEED"
$ws.Range("N2").Value = "SPEED"

# Update the selected/active cell to reflect the authored state (O17)
$ws.Range("O17").Select()
